$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fit-parameter cells (I3, I4, I5)
$ws.Range("I3").Value = 27099700000
$ws.Range("I4").Value = 1.53691
$ws.Range("I5").Value = -1.36972

# Clear the number-format style previously applied to I4/I5 (now plain numbers)
$ws.Range("I4").Style = "Normal"
$ws.Range("I5").Style = "Normal"

# Update the formula in E5 and extend it down to E6 and E7
$ws.Range("E5").Formula = "=`$I`$3*Table1[[#This Row],[Mean R]]^-`$I`$4 + `$I`$5"
$ws.Range("E6").Formula = "=`$I`$3*Table1[[#This Row],[Mean R]]^-`$I`$4 + `$I`$5"
$ws.Range("E7").Formula = "=`$I`$3*Table1[[#This Row],[Mean R]]^-`$I`$4 + `$I`$5"

# Update the selection to J14
$ws.Range("J14").Select()
